$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Reorder columns D/E/F (Ship, Faction, Unique -> Unique, Ship, Faction) ---
# Capture current values before overwriting them.
$ship = @{}
$faction = @{}
$unique = @{}
for ($r = 2; $r -le 4; $r++) {
  $ship[$r] = $ws.Cells.Item($r, 4).Value2
  $faction[$r] = $ws.Cells.Item($r, 5).Value2
  $unique[$r] = $ws.Cells.Item($r, 6).Value2
}

# New headers for D/E/F.
$ws.Cells.Item(1, 4).Value = "Unique"
$ws.Cells.Item(1, 5).Value = "Ship"
$ws.Cells.Item(1, 6).Value = "Faction"

# Write back data rows in the new column order.
for ($r = 2; $r -le 4; $r++) {
  $ws.Cells.Item($r, 4).Value = $unique[$r]
  $ws.Cells.Item($r, 5).Value = $ship[$r]
  $ws.Cells.Item($r, 6).Value = $faction[$r]
}

# --- Remove the old "Limited" column (G) entirely ---
$ws.Range("G1:G4").Delete()

# --- Fill in the missing pilot ability text for row 4 ---
$ws.Cells.Item(4, 2).Value = "He's Just a Guy"

# --- Fix up the view: select B4 (also clears the stale topLeftCell) ---
$ws.Range("B4").Select()
